# Append the daily block of hospital records for 2020-06-06 (Excel serial 43988),
# replicating the prior day's 20-row block (rows 1253:1272, date 43987) down into
# rows 1273:1292 and then updating the date and the few occupancy counts that changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the previous day's 20-row block (with all formatting/styles) into the new rows.
$src = $ws.Range("A1253:H1272")
$dst = $ws.Range("A1273:H1292")
$src.Copy($dst)

# New date for the appended block: 2020-06-06 -> serial 43988
$ws.Range("A1273:A1292").Value = 43988

# Update the handful of "camas_ocupadas_total" (column C) values that changed vs. the
# previous day's figures.
$ws.Cells.Item(1273, 3).Value = 8
$ws.Cells.Item(1274, 3).Value = 33
$ws.Cells.Item(1278, 3).Value = 5
$ws.Cells.Item(1280, 3).Value = 9
$ws.Cells.Item(1281, 3).Value = 4
$ws.Cells.Item(1284, 3).Value = 4
$ws.Cells.Item(1288, 3).Value = 1
$ws.Cells.Item(1292, 3).Value = 1
